$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 247.57143
$ws.Range("I11").Value = 247.57143
$ws.Range("K11").Value = 247.57143
$ws.Range("M11").Value = -107.57143

$ws.Range("H40").Value = 911353.25
$ws.Range("I40").Value = 911353.25
$ws.Range("K40").Value = 911353.25
$ws.Range("M40").Value = -911178.25

$ws.Range("H41").Value = 294.16666
$ws.Range("I41").Value = 533
$ws.Range("J41").Value = 55.333332
$ws.Range("K41").Value = 533
$ws.Range("L41").Value = 55.333332
$ws.Range("M41").Value = -93
$ws.Range("N41").Value = -935.333332

$ws.Range("H64").Value = 11128.429
$ws.Range("I64").Value = 12199.8
$ws.Range("K64").Value = 12199.8
$ws.Range("M64").Value = -11951.8

$ws.Range("H67").Value = 11128.429
$ws.Range("I67").Value = 12199.8
$ws.Range("K67").Value = 12199.8
$ws.Range("M67").Value = -11341.8

$ws.Range("H96").Value = 2320.8667
$ws.Range("I96").Value = 2002
$ws.Range("J96").Value = 2685.2856
$ws.Range("K96").Value = 6006
$ws.Range("L96").Value = 8055.8568
$ws.Range("M96").Value = -4633
$ws.Range("N96").Value = -10801.8568

$ws.Range("H113").Value = 5058.5
$ws.Range("J113").Value = 5712.4287
$ws.Range("L113").Value = 5712.4287
$ws.Range("N113").Value = -12220.4287

$ws.Range("H116").Value = 2904.3076
$ws.Range("I116").Value = 3400
$ws.Range("J116").Value = 2326
$ws.Range("K116").Value = 3400
$ws.Range("L116").Value = 2326
$ws.Range("M116").Value = 42
$ws.Range("N116").Value = -9210

$ws.Range("H132").Value = 2496.639
$ws.Range("I132").Value = 2051.8333
$ws.Range("K132").Value = 6155.499899999999
$ws.Range("M132").Value = -3625.499899999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3584.745
$ws.Range("I32").Value = 2745.2126
$ws.Range("J32").Value = 13449.25
$ws.Range("K32").Value = 2745.2126
$ws.Range("L32").Value = 13449.25
$ws.Range("M32").Value = -2458.2126
$ws.Range("N32").Value = -14023.25

$ws.Range("H45").Value = 2072.1428
$ws.Range("J45").Value = 2469
$ws.Range("L45").Value = 2469
$ws.Range("N45").Value = -3223

$ws.Range("H74").Value = 52691450
$ws.Range("I74").Value = 62570908
$ws.Range("J74").Value = 1003
$ws.Range("K74").Value = 62570908
$ws.Range("L74").Value = 1003
$ws.Range("M74").Value = -62570034
$ws.Range("N74").Value = -2751

$ws.Range("H77").Value = 52691450
$ws.Range("I77").Value = 62570908
$ws.Range("J77").Value = 1003
$ws.Range("K77").Value = 312854540
$ws.Range("L77").Value = 5015
$ws.Range("M77").Value = -312850172
$ws.Range("N77").Value = -13751

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 13004
$ws.Range("I82").Value = 13004
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 13004
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -12621
$ws.Range("N82").ClearContents()

$ws.Range("H85").Value = 13004
$ws.Range("I85").Value = 13004
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 13004
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -11678
$ws.Range("N85").ClearContents()

$ws.Range("H86").Value = 13516.375
$ws.Range("I86").Value = 5790.857
$ws.Range("J86").Value = 67595
$ws.Range("K86").Value = 5790.857
$ws.Range("L86").Value = 67595
$ws.Range("M86").Value = -4667.857
$ws.Range("N86").Value = -69841

$ws.Range("H89").Value = 13516.375
$ws.Range("I89").Value = 5790.857
$ws.Range("J89").Value = 67595
$ws.Range("K89").Value = 28954.285
$ws.Range("L89").Value = 337975
$ws.Range("M89").Value = -23338.285
$ws.Range("N89").Value = -349207

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("I7").Value = 112.5
$ws.Range("J7").Value = 322.06668
$ws.Range("K7").Value = 112.5
$ws.Range("L7").Value = 322.06668
$ws.Range("M7").Value = 0.5
$ws.Range("N7").Value = -548.06668

$ws.Range("H22").Value = 21580.2
$ws.Range("I22").Value = 33900.332
$ws.Range("J22").Value = 3100
$ws.Range("K22").Value = 33900.332
$ws.Range("L22").Value = 3100
$ws.Range("M22").Value = -33550.332
$ws.Range("N22").Value = -3800

$ws.Range("H98").Value = 36153
$ws.Range("J98").Value = 36153
$ws.Range("L98").Value = 36153
$ws.Range("N98").Value = -40645

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 266.16666
$ws.Range("I7").Value = 302.10526
$ws.Range("J7").Value = 129.6
$ws.Range("K7").Value = 906.3157799999999
$ws.Range("L7").Value = 388.8
$ws.Range("M7").Value = -794.3157799999999
$ws.Range("N7").Value = -612.8

$ws.Range("H120").Value = 11500
$ws.Range("I120").Value = 11500
$ws.Range("K120").Value = 34500
$ws.Range("M120").Value = -29662

$ws.Range("H129").Value = 4506.5557
$ws.Range("I129").Value = 5011.143
$ws.Range("K129").Value = 15033.429
$ws.Range("M129").Value = -10033.429

$ws.Range("H131").Value = 24958.875
$ws.Range("J131").Value = 4527.5
$ws.Range("L131").Value = 13582.5
$ws.Range("N131").Value = -23662.5

$ws.Range("H140").Value = 2398.7
$ws.Range("I140").Value = 2240.125
$ws.Range("K140").Value = 6720.375
$ws.Range("M140").Value = -1540.375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7581.273
$ws.Range("I70").Value = 6965.933
$ws.Range("K70").Value = 6965.933
$ws.Range("M70").Value = -6695.933

$ws.Range("H73").Value = 7581.273
$ws.Range("I73").Value = 6965.933
$ws.Range("K73").Value = 6965.933
$ws.Range("M73").Value = -6029.933

$ws.Range("H80").Value = 4336.2354
$ws.Range("I80").Value = 4423
$ws.Range("J80").Value = 4128
$ws.Range("K80").Value = 4423
$ws.Range("L80").Value = 4128
$ws.Range("M80").Value = -3425
$ws.Range("N80").Value = -6124

$ws.Range("H83").Value = 4336.2354
$ws.Range("I83").Value = 4423
$ws.Range("J83").Value = 4128
$ws.Range("K83").Value = 22115
$ws.Range("L83").Value = 20640
$ws.Range("M83").Value = -17123
$ws.Range("N83").Value = -30624

$ws.Range("H109").Value = 33408.168
$ws.Range("J109").Value = 32987.25
$ws.Range("L109").Value = 32987.25
$ws.Range("N109").Value = -35067.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 535.44446
$ws.Range("I16").Value = 452.375
$ws.Range("K16").Value = 452.375
$ws.Range("M16").Value = -282.375

$ws.Range("H22").Value = 2524
$ws.Range("I22").Value = 1510
$ws.Range("J22").Value = 2929.6
$ws.Range("K22").Value = 1510
$ws.Range("L22").Value = 2929.6
$ws.Range("M22").Value = -1215
$ws.Range("N22").Value = -3519.6

$ws.Range("H27").Value = 2524
$ws.Range("I27").Value = 1510
$ws.Range("J27").Value = 2929.6
$ws.Range("K27").Value = 1510
$ws.Range("L27").Value = 2929.6
$ws.Range("M27").Value = -1403
$ws.Range("N27").Value = -3143.6

$ws.Range("H46").Value = 844.0540999999999
$ws.Range("I46").Value = 516.76666
$ws.Range("K46").Value = 516.76666
$ws.Range("M46").Value = -328.76666

$ws.Range("H51").Value = 40495
$ws.Range("J51").Value = 40495
$ws.Range("L51").Value = 40495
$ws.Range("N51").Value = -41451

$ws.Range("H68").Value = 3409.1667
$ws.Range("I68").Value = 2360
$ws.Range("J68").Value = 4458.3335
$ws.Range("K68").Value = 2360
$ws.Range("L68").Value = 4458.3335
$ws.Range("M68").Value = -1611
$ws.Range("N68").Value = -5956.3335

$ws.Range("H71").Value = 3409.1667
$ws.Range("I71").Value = 2360
$ws.Range("J71").Value = 4458.3335
$ws.Range("K71").Value = 11800
$ws.Range("L71").Value = 22291.6675
$ws.Range("M71").Value = -8056
$ws.Range("N71").Value = -29779.6675

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 5005
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 5005
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 5005
$ws.Range("N7").Value = -5231
$ws.Range("M7").ClearContents()

$ws.Range("H13").Value = 2415.5
$ws.Range("I13").Value = 2623.25
$ws.Range("K13").Value = 2623.25
$ws.Range("M13").Value = -2483.25

$ws.Range("H22").Value = 1932.3334
$ws.Range("I22").Value = 1800
$ws.Range("J22").Value = 1998.5
$ws.Range("K22").Value = 1800
$ws.Range("L22").Value = 1998.5
$ws.Range("M22").Value = -1507
$ws.Range("N22").Value = -2584.5

$ws.Range("H26").Value = 2000
$ws.Range("J26").Value = 2000
$ws.Range("L26").Value = 2000
$ws.Range("N26").Value = -2586

$ws.Range("H29").Value = 2000
$ws.Range("J29").Value = 1000
$ws.Range("L29").Value = 1000
$ws.Range("N29").Value = -1580

$ws.Range("H43").Value = 29050
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 29050
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 29050
$ws.Range("N43").Value = -29348
$ws.Range("M43").ClearContents()

$ws.Range("H135").Value = 100000000
$ws.Range("J135").Value = 100000000
$ws.Range("L135").Value = 100000000
$ws.Range("N135").Value = -100010140

$ws.Range("H140").Value = 71808
$ws.Range("J140").Value = 71808
$ws.Range("L140").Value = 71808
$ws.Range("N140").Value = -82168

$ws.Range("H141").Value = 64913.75
$ws.Range("J141").Value = 25715
$ws.Range("L141").Value = 25715
$ws.Range("N141").Value = -36075
